# Update the NATMI LR-pair table (Il4-Cd53) with the new TPM-based results.
# The re-run of the pipeline dropped "MuSCs" as a possible Target cluster,
# so the rows where column D = "MuSCs" are removed entirely, and the
# remaining rows' statistics (columns E:T) are recomputed against the
# smaller background population.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four rows whose Target cluster (column D) is "MuSCs".
# Delete from the bottom up so earlier row numbers stay valid.
$ws.Rows("12").Delete()
$ws.Rows("9").Delete()
$ws.Rows("6").Delete()
$ws.Rows("3").Delete()

# New recalculated values (columns E through T) for the 8 remaining data
# rows (now rows 2-9):
#   row 2: ECs   / ECs
#   row 3: ECs   / Resolving-Mac
#   row 4: FAPs  / ECs
#   row 5: FAPs  / Resolving-Mac
#   row 6: MuSCs / ECs
#   row 7: MuSCs / Resolving-Mac
#   row 8: Resolving-Mac / ECs
#   row 9: Resolving-Mac / Resolving-Mac
$rowsData = @(
    @(3,1,0.9458723333333334,2.837617,0.1873686327665471,0.1873686327665471,1,0.3333333333333333,0.096592,0.289776,0.0005851292956313063,0.0005851292956313063,0.09136370042133334,0.822273303792,0.0001096348761140906,0.0001096348761140906),
    @(3,1,0.9458723333333334,2.837617,0.1873686327665471,0.1873686327665471,3,1,164.9814526666667,494.944358,0.9994148707043687,0.9994148707043686,156.0513915905429,1404.462524314886,0.187258997890433,0.187258997890433),
    @(3,1,1.964334,5.893002,0.3891165466060174,0.3891165466060174,1,0.3333333333333333,0.096592,0.289776,0.0005851292956313063,0.0005851292956313063,0.189738949728,1.707650547552,0.0002276834908340653,0.0002276834908340653),
    @(3,1,1.964334,5.893002,0.3891165466060174,0.3891165466060174,3,1,164.9814526666667,494.944358,0.9994148707043687,0.9994148707043686,324.078676842524,2916.708091582716,0.3888888631151833,0.3888888631151833),
    @(3,1,1.499502333333333,4.498507,0.2970376573303378,0.2970376573303378,1,0.3333333333333333,0.096592,0.289776,0.0005851292956313063,0.0005851292956313063,0.1448399293813333,1.303559364432,0.0001738054352096739,0.0001738054352096739),
    @(3,1,1.499502333333333,4.498507,0.2970376573303378,0.2970376573303378,3,1,164.9814526666667,494.944358,0.9994148707043687,0.9994148707043686,247.3900732303896,2226.510659073506,0.2968638518951282,0.2968638518951282),
    @(3,1,0.6384806666666667,1.915442,0.1264771632970977,0.1264771632970977,1,0.3333333333333333,0.096592,0.289776,0.0005851292956313063,0.0005851292956313063,0.06167212455466667,0.555049120992,0.00007400549347347647,0.00007400549347347647),
    @(3,1,0.6384806666666667,1.915442,0.1264771632970977,0.1264771632970977,3,1,164.9814526666667,494.944358,0.9994148707043687,0.9994148707043686,105.3374678862485,948.0372109762361,0.1264031578036242,0.1264031578036242)
)

$r = 2
foreach ($rowData in $rowsData) {
    $c = 5
    foreach ($val in $rowData) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
